$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.197.88"
$ws.Range("E2").Value = "  +0.31%  "

$ws.Range("D3").Value = "2.559.02"
$ws.Range("E3").Value = "  +0.71%  "

$ws.Range("E4").Value = "  +0.07%  "

$c = $ws.Range("D5")
$c.Formula = "'583.98"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.55%  "

$c = $ws.Range("D6")
$c.Formula = "'148.00"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.20%  "

$ws.Range("E7").Value = "  +0.06%  "

$c = $ws.Range("D8")
$c.Formula = "'0.586"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.09%  "

$c = $ws.Range("D9")
$c.Formula = "'0.109"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +3.79%  "

$c = $ws.Range("D10")
$c.Formula = "'5.62"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.01%  "

$ws.Range("E11").Value = "  +0.25%  "

$ws.Range("E12").Value = "  +0.21%  "

$c = $ws.Range("D13")
$c.Formula = "'27.64"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.56%  "

$ws.Range("D14").Value = "3.016.45"
$ws.Range("E14").Value = "  +0.78%  "

$ws.Range("D15").Value = "63.050.59"
$ws.Range("E15").Value = "  +0.20%  "

$c = $ws.Range("D16")
$c.Formula = "'0.0000148"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +4.46%  "

$ws.Range("D17").Value = "2.556.05"
$ws.Range("E17").Value = "  +0.43%  "

$ws.Range("E18").Value = "  -1.01%  "

$c = $ws.Range("D19")
$c.Formula = "'342.40"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +2.68%  "

$c = $ws.Range("D20")
$c.Formula = "'4.43"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +3.18%  "

$ws.Range("E21").Value = "  +1.17%  "

$ws.Range("E22").Value = "  -0.13%  "

$c = $ws.Range("D23")
$c.Formula = "'66.51"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +2.64%  "

$ws.Range("E24").Value = "  +2.81%  "

$ws.Range("D25").Value = "2.685.76"
$ws.Range("E25").Value = "  +0.78%  "

$c = $ws.Range("D26")
$c.Formula = "'0.171"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +1.29%  "

$c = $ws.Range("D27")
$c.Formula = "'8.13"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +13.28%  "

$c = $ws.Range("D28")
$c.Formula = "'8.55"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +1.89%  "

$ws.Range("E30").Value = "  +0.56%  "

$c = $ws.Range("D31")
$c.Formula = "'1.99"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +8.05%  "

$ws.Range("D32").Value = "0.0₃0825"
$ws.Range("E32").Value = "  +1.37%  "

$c = $ws.Range("D33")
$c.Formula = "'176.95"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.22%  "

$c = $ws.Range("D34")
$c.Formula = "'440.86"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +6.41%  "

$ws.Range("E35").Value = "  +2.26%  "

$ws.Range("E36").Value = "  +2.27%  "

$c = $ws.Range("D37")
$c.Formula = "'19.27"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +2.28%  "

$ws.Range("E38").Value = "  +3.41%  "

$ws.Range("E40").Value = "  +0.10%  "

$ws.Range("E41").Value = "  +0.12%  "

$c = $ws.Range("D42")
$c.Formula = "'150.81"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.37%  "

$c = $ws.Range("D43")
$c.Formula = "'3.84"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +2.53%  "

$ws.Range("E44").Value = "  +2.39%  "

$c = $ws.Range("D45")
$c.Formula = "'0.0549"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +6.10%  "

$ws.Range("E46").Value = "  +1.38%  "

$ws.Range("E47").Value = "  +1.30%  "

$ws.Range("E48").Value = "  +2.74%  "

$c = $ws.Range("D49")
$c.Formula = "'18.42"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.65%  "

$ws.Range("E50").Value = "  -2.08%  "

$ws.Range("E51").Value = "  -0.29%  "
